# Updates the "cryptos" price/volume table to the latest scrape.
# Source data: coinranking.com rankings 0-49 (rows 2-51), columns:
#   B = Coin name, C = Link, D = Price (text, locale-formatted), E = Volume(1h) (text, padded %).
#
# Two pairs of rows had their rankings swap places versus the prior run
# (Bittensor/Aptos at rows 31/32, Stellar/VeChain at rows 48/49), so those
# rows get full B:E replacements; every other touched row only gets new
# D/E figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many "Price" values look numeric ("1.00", "0.0960", "146.00") but the
# column stores them as plain text (leading/trailing zeros, locale dot-
# grouping like "62.985.80" must survive verbatim). Assigning such a
# string straight to .Value lets Excel's input-parsing coerce it into a
# real number and drop the formatting, so those are written with a
# leading apostrophe (forces text entry) and then have their style reset
# to match an untouched "Price" cell, undoing the quote-prefix style Excel
# applies automatically - the net result is plain text content under the
# same default "General" style the cell started with.
$plainPriceStyle = $ws.Range("D27").Style

function Set-CellText($Cell, $Text, $ForceText) {
    $range = $ws.Range($Cell)
    if ($ForceText) {
        $range.Value = "'" + $Text
        $range.Style = $plainPriceStyle
    } else {
        $range.Value = $Text
    }
}

$updates = @(
    @{ Cell = "D2"; Value = "62.985.80"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -2.33%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "2.621.61"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -2.19%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  -0.02%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "603.72"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +1.07%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "146.00"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -1.61%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  -0.01%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  -1.52%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "2.621.03"; ForceText = $false },
    @{ Cell = "E9"; Value = "  -2.24%  "; ForceText = $false },
    @{ Cell = "E10"; Value = "  -0.57%  "; ForceText = $false },
    @{ Cell = "E11"; Value = "  -1.35%  "; ForceText = $false },
    @{ Cell = "E12"; Value = "  +0.03%  "; ForceText = $false },
    @{ Cell = "E13"; Value = "  +0.49%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "27.15"; ForceText = $true },
    @{ Cell = "E14"; Value = "  -3.13%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "3.091.73"; ForceText = $false },
    @{ Cell = "E15"; Value = "  -2.28%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "62.885.54"; ForceText = $false },
    @{ Cell = "E16"; Value = "  -2.36%  "; ForceText = $false },
    @{ Cell = "E17"; Value = "  -2.41%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "2.636.25"; ForceText = $false },
    @{ Cell = "E18"; Value = "  -2.26%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "11.27"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -1.66%  "; ForceText = $false },
    @{ Cell = "E20"; Value = "  +0.92%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "339.76"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -2.17%  "; ForceText = $false },
    @{ Cell = "E22"; Value = "  -1.05%  "; ForceText = $false },
    @{ Cell = "E23"; Value = "  -0.09%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "5.56"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -4.79%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "66.52"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -3.62%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "1.61"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -3.71%  "; ForceText = $false },
    @{ Cell = "E27"; Value = "  -4.84%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "8.65"; ForceText = $true },
    @{ Cell = "E28"; Value = "  +1.25%  "; ForceText = $false },
    @{ Cell = "E29"; Value = "  -2.66%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E30"; Value = "  +0.25%  "; ForceText = $false },
    @{ Cell = "B31"; Value = "Aptos"; ForceText = $false },
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; ForceText = $false },
    @{ Cell = "D31"; Value = "7.91"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -1.23%  "; ForceText = $false },
    @{ Cell = "B32"; Value = "Bittensor"; ForceText = $false },
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; ForceText = $false },
    @{ Cell = "D32"; Value = "535.82"; ForceText = $true },
    @{ Cell = "E32"; Value = "  +0.62%  "; ForceText = $false },
    @{ Cell = "E33"; Value = "  +1.13%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "1.74"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -2.44%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "0.0₃0802"; ForceText = $false },
    @{ Cell = "E35"; Value = "  -3.18%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "5.18"; ForceText = $true },
    @{ Cell = "E36"; Value = "  +10.22%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "168.66"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -3.94%  "; ForceText = $false },
    @{ Cell = "E38"; Value = "  -0.09%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "0.402"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -0.06%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "19.00"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -1.98%  "; ForceText = $false },
    @{ Cell = "E41"; Value = "  +5.00%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  -0.04%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "169.51"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -2.82%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "3.74"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -1.59%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "22.35"; ForceText = $true },
    @{ Cell = "E45"; Value = "  +1.45%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "0.0567"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +2.72%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.623"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -2.28%  "; ForceText = $false },
    @{ Cell = "B48"; Value = "VeChain"; ForceText = $false },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; ForceText = $false },
    @{ Cell = "D48"; Value = "0.0240"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -1.05%  "; ForceText = $false },
    @{ Cell = "B49"; Value = "Stellar"; ForceText = $false },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; ForceText = $false },
    @{ Cell = "D49"; Value = "0.0960"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -0.61%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "18.46"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -2.66%  "; ForceText = $false },
    @{ Cell = "E51"; Value = "  -0.38%  "; ForceText = $false }
)

foreach ($u in $updates) {
    Set-CellText $u.Cell $u.Value $u.ForceText
}
